$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J4").Value = 39
$ws.Range("J5").Value = 20
$ws.Range("J6").Value = 5
$ws.Range("J7").Value = 10
$ws.Range("J8").Value = 40
$ws.Range("J9").Value = 15
$ws.Range("J10").Value = 129

$ws.Range("J14").Value = 30
$ws.Range("J15").Value = 30
$ws.Range("J16").Value = 30
$ws.Range("J17").Value = 60
$ws.Range("J18").Value = 15
$ws.Range("J19").Value = 45
$ws.Range("J20").Value = 30
$ws.Range("J21").Value = 30
$ws.Range("J22").Value = 30
$ws.Range("J23").Value = 45

$ws.Range("J31").Value = 2
$ws.Range("J45").Value = 2
